# remove rubber from BoM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Rubber" / "to compress to bike rack" / "to compress between bike
# rack and enclosure. We used ninjaflex 3D print filament") is being removed
# entirely from the Bill of Materials. Clearing every cell in that row drops
# the now-empty <row> element from sheetData (rows below it, e.g. row 17,
# keep their original row numbers).
$ws.Range("A12:E12").ClearContents()

# The author's selection ended up on C14 after the edit.
$ws.Range("C14").Select()
